$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.386222362518311
$ws.Range("B1").Value = 1.745376229286194
$ws.Range("C1").Value = 6.746384143829346
$ws.Range("D1").Value = 1.623165965080261
$ws.Range("E1").Value = 0.9638630747795105
